# Apply the "collection feature" edits to constant_vocab_mapping.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant_vocab_mapping")

# A17: rename the vocab term from og_group_ref -> field_collection_field
$ws.Range("A17").Value = "field_collection_field"

# B17: update the id value
$ws.Range("B17").Value = 1566

# Update the active selection to B17 (was D18)
$ws.Range("B17").Select()
